# Regenerate merged AHB files
# - rename header columns from *_old / *_new to *_FV2410 / *_FV2504
# - wrap the data range in an Excel Table (Table1)
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newHeaders = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410",
    "diff",
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $newHeaders[$i]
}

# Turn the header+data range into a native Excel Table
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U55"), $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (pane split below row 1)
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Host "edit applied"
